{"js": "// Replace each two-digit-division formula with its updated value.\n// Old/new text pairs taken from the authoritative diff (order-independent,\n// each old value is unique in the document so search+replace is unambiguous).\nconst replacements = [\n  [\"18\u00f79=2, 0\", \"50\u00f74=12, 2\"],\n  [\"91\u00f75=18, 1\", \"99\u00f77=14, 1\"],\n  [\"14\u00f75=2, 4\", \"85\u00f77=12, 1\"],\n  [\"84\u00f78=10, 4\", \"47\u00f74=11, 3\"],\n  [\"53\u00f79=5, 8\", \"66\u00f78=8, 2\"],\n  [\"65\u00f79=7, 2\", \"86\u00f75=17, 1\"],\n  [\"23\u00f77=3, 2\", \"28\u00f75=5, 3\"],\n  [\"30\u00f76=5, 0\", \"74\u00f76=12, 2\"],\n  [\"83\u00f76=13, 5\", \"19\u00f78=2, 3\"],\n  [\"44\u00f74=11, 0\", \"52\u00f73=17, 1\"],\n  [\"13\u00f73=4, 1\", \"31\u00f74=7, 3\"],\n  [\"42\u00f76=7, 0\", \"67\u00f78=8, 3\"],\n  [\"25\u00f78=3, 1\", \"59\u00f79=6, 5\"],\n  [\"76\u00f75=15, 1\", \"90\u00f78=11, 2\"],\n  [\"94\u00f79=10, 4\", \"77\u00f72=38, 1\"],\n  [\"87\u00f74=21, 3\", \"57\u00f78=7, 1\"],\n  [\"31\u00f78=3, 7\", \"38\u00f75=7, 3\"],\n  [\"97\u00f77=13, 6\", \"58\u00f73=19, 1\"],\n  [\"84\u00f72=42, 0\", \"49\u00f72=24, 1\"],\n  [\"23\u00f75=4, 3\", \"99\u00f74=24, 3\"],\n  [\"96\u00f73=32, 0\", \"83\u00f72=41, 1\"],\n  [\"13\u00f79=1, 4\", \"35\u00f74=8, 3\"],\n  [\"64\u00f76=10, 4\", \"98\u00f75=19, 3\"],\n  [\"96\u00f77=13, 5\", \"29\u00f77=4, 1\"],\n  [\"66\u00f79=7, 3\", \"38\u00f73=12, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-division formula with its updated value.\n# Old/new text pairs taken from the authoritative diff (order-independent,\n# each old value is unique in the document so Find/Replace is unambiguous).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('18\u00f79=2, 0', '50\u00f74=12, 2'),\n    @('91\u00f75=18, 1', '99\u00f77=14, 1'),\n    @('14\u00f75=2, 4', '85\u00f77=12, 1'),\n    @('84\u00f78=10, 4', '47\u00f74=11, 3'),\n    @('53\u00f79=5, 8', '66\u00f78=8, 2'),\n    @('65\u00f79=7, 2', '86\u00f75=17, 1'),\n    @('23\u00f77=3, 2', '28\u00f75=5, 3'),\n    @('30\u00f76=5, 0', '74\u00f76=12, 2'),\n    @('83\u00f76=13, 5', '19\u00f78=2, 3'),\n    @('44\u00f74=11, 0', '52\u00f73=17, 1'),\n    @('13\u00f73=4, 1', '31\u00f74=7, 3'),\n    @('42\u00f76=7, 0', '67\u00f78=8, 3'),\n    @('25\u00f78=3, 1', '59\u00f79=6, 5'),\n    @('76\u00f75=15, 1', '90\u00f78=11, 2'),\n    @('94\u00f79=10, 4', '77\u00f72=38, 1'),\n    @('87\u00f74=21, 3', '57\u00f78=7, 1'),\n    @('31\u00f78=3, 7', '38\u00f75=7, 3'),\n    @('97\u00f77=13, 6', '58\u00f73=19, 1'),\n    @('84\u00f72=42, 0', '49\u00f72=24, 1'),\n    @('23\u00f75=4, 3', '99\u00f74=24, 3'),\n    @('96\u00f73=32, 0', '83\u00f72=41, 1'),\n    @('13\u00f79=1, 4', '35\u00f74=8, 3'),\n    @('64\u00f76=10, 4', '98\u00f75=19, 3'),\n    @('96\u00f77=13, 5', '29\u00f77=4, 1'),\n    @('66\u00f79=7, 3', '38\u00f73=12, 2'),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute(\n        $oldText,   # FindText\n        $false,     # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $newText,   # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n"}
